# Replace the four "Perseus havainnointijaksot vuonna 2018: ..." date-range
# paragraphs with the single, plain-formatted Gemini date string, collapsing
# every run in each paragraph (including the trailing "Ennen kuin..." text
# and the www.globeatnight.org/finding hyperlink in the third occurrence)
# down to one run with no explicit run formatting.

$d = $word.ActiveDocument
$newText = "havainnointijaksot vuonna Gemini: 14.-23. Helmikuuta 14.-24."

# Locate every paragraph whose text starts with the old Perseus/2018 blurb.
# (Word's Paragraphs collection is 1-based and also walks table/other
# story paragraphs, so we search rather than hard-code indices.)
$count = $d.Paragraphs.Count
$targets = @()
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Perseus havainnointijaksot vuonna 2018*") {
        $targets += $i
    }
}

foreach ($idx in $targets) {
    $p = $d.Paragraphs($idx)
    $pr = $p.Range
    # Re-materialize the range from plain offsets: setting .Text on a
    # Range obtained straight from Paragraphs(n).Range only touches the
    # first run, so rebuild an equivalent Range via the document instead.
    $start = $pr.Start
    $end = $pr.End
    $r = $d.Range($start, $end)

    # Wipe every run (text + hyperlink) inside the paragraph while
    # preserving the paragraph mark / pPr, then insert one fresh run
    # holding the new text with default (inherited) run formatting.
    $r.Text = ""
    $r.InsertAfter($newText)
}
